$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Default_Values")

# Update default parameter values (row numbers correspond to parameter rows on the sheet)
$ws.Range("B3").Value = 99999       # AnnualEmissionLimit
$ws.Range("B17").Value = 0.05       # DiscountRate
$ws.Range("B23").Value = 99999      # ModelPeriodEmissionLimit
$ws.Range("B43").Value = 99999      # TotalAnnualMaxCapacity
$ws.Range("B44").Value = 99999      # TotalAnnualMaxCapacityInvestment
$ws.Range("B48").Value = 99999      # TotalTechnologyAnnualActivityUpperLimit
$ws.Range("B50").Value = 99999      # TotalTechnologyModelPeriodActivityUpperLimit
